$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet name / title to reflect the new "through" date
$wb.Worksheets.Item(1).Name = "Through 2022-12-11"

# Update the row label for December
$ws.Range("A13").Value = "December (through 12-11)"

# Update December row (row 13) values
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 39
$ws.Range("F13").Value = 16
$ws.Range("G13").Value = 55
$ws.Range("H13").Value = 85
$ws.Range("I13").Value = 50

# Update Total row (row 14) values
$ws.Range("B14").Value = 301
$ws.Range("C14").Value = 593
$ws.Range("D14").Value = 860
$ws.Range("F14").Value = 550
$ws.Range("G14").Value = 1319
$ws.Range("H14").Value = 1728
$ws.Range("I14").Value = 1566
